# Actualiza notas estudiantes taller 3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Garcia Lopez, Jose Manuel (fila 10)
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 9.8
$ws.Range("D10").Value = 9

# Garcia Vargas, Juan Francesco (fila 11)
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = 7
$ws.Range("D11").Value = 0
# Nota inconsistente / en cero -> resaltar en rojo para revisar de nuevo
$ws.Range("D11").Font.Color = 255

# Gonzalez Castrillon, Miguel Angel (fila 13)
$ws.Range("B13").Value = 10
$ws.Range("C13").Value = 10
$ws.Range("D13").Value = 10

# Rojas Mejia, Juan Miguel (fila 18)
$ws.Range("B18").Value = 10
$ws.Range("C18").Value = 10
$ws.Range("D18").Value = 10

# Toro Trujillo, Juan Jose (fila 20)
$ws.Range("B20").Value = 10
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 10

# Actualiza la celda seleccionada / vista de la hoja
$ws.Range("B9").Select()
